# Generate Report for Handoff
# Adds two new handed-off files (a5161b99-417d-4b7d-9d2e-70884dde2daa.md and
# a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md) as rows 4/5 on every sheet, extends
# the three ListObjects to cover the new rows, and wires up the matching
# hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A4").Value = "a5161b99-417d-4b7d-9d2e-70884dde2daa.md"
$ov.Range("B4").Value = "e2e\a5161b99-417d-4b7d-9d2e-70884dde2daa.md"
$ov.Range("C4").Value = ".md"
$ov.Range("E4").Value = "Ready for handoff"
$ov.Range("F4").Value = "Ready for handoff"
$ov.Range("G4").Value = "2016-08-13 20:52:54"
$ov.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Range("A5").Value = "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md"
$ov.Range("B5").Value = "e2e\a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md"
$ov.Range("C5").Value = ".md"
$ov.Range("E5").Value = "Ready for handoff"
$ov.Range("F5").Value = "Ready for handoff"
$ov.Range("G5").Value = "2016-08-13 20:52:54"
$ov.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a5161b99417d4b7d9d2e70884dde2daa0000000/e2e/a5161b99-417d-4b7d-9d2e-70884dde2daa.md", "", "", "e2e\a5161b99-417d-4b7d-9d2e-70884dde2daa.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a7a29e68ec614dc3a4d13cbe7cc940640000000/e2e/a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md", "", "", "e2e\a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md") | Out-Null

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A4").Value = "a5161b99-417d-4b7d-9d2e-70884dde2daa.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "e2e"
$zh.Range("E4").Value = "ht"
$zh.Range("F4").Value = "'False"
$zh.Range("G4").Value = "a5161b99-417d-4b7d-9d2e-70884dde2daa.05dabd5d46a2cdb06abd1ac3e6b604d127148806.zh-cn.xlf"
$zh.Range("H4").Value = "2016-08-13 20:52:46"
$zh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("K4").Value = "0001-01-01 00:00:00"
$zh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("M4").Value = "'True"
$zh.Range("O4").Value = "'False"

$zh.Range("A5").Value = "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "e2e"
$zh.Range("E5").Value = "ht"
$zh.Range("F5").Value = "'False"
$zh.Range("G5").Value = "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.caaeabc7d9e86c70c4ab018d5cbde5253e51c6fa.zh-cn.xlf"
$zh.Range("H5").Value = "2016-08-13 20:52:46"
$zh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("K5").Value = "0001-01-01 00:00:00"
$zh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("M5").Value = "'True"
$zh.Range("O5").Value = "'False"

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a5161b99417d4b7d9d2e70884dde2daa0000000/e2e/a5161b99-417d-4b7d-9d2e-70884dde2daa.md", "", "", "a5161b99-417d-4b7d-9d2e-70884dde2daa.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a7a29e68ec614dc3a4d13cbe7cc940640000000/e2e/a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md", "", "", "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md") | Out-Null

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A4").Value = "a5161b99-417d-4b7d-9d2e-70884dde2daa.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "e2e"
$de.Range("E4").Value = "ht"
$de.Range("F4").Value = "'False"
$de.Range("G4").Value = "a5161b99-417d-4b7d-9d2e-70884dde2daa.05dabd5d46a2cdb06abd1ac3e6b604d127148806.de-de.xlf"
$de.Range("H4").Value = "2016-08-13 20:52:54"
$de.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("K4").Value = "0001-01-01 00:00:00"
$de.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("M4").Value = "'True"
$de.Range("O4").Value = "'False"

$de.Range("A5").Value = "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "e2e"
$de.Range("E5").Value = "ht"
$de.Range("F5").Value = "'False"
$de.Range("G5").Value = "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.caaeabc7d9e86c70c4ab018d5cbde5253e51c6fa.de-de.xlf"
$de.Range("H5").Value = "2016-08-13 20:52:54"
$de.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("K5").Value = "0001-01-01 00:00:00"
$de.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("M5").Value = "'True"
$de.Range("O5").Value = "'False"

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a5161b99417d4b7d9d2e70884dde2daa0000000/e2e/a5161b99-417d-4b7d-9d2e-70884dde2daa.md", "", "", "a5161b99-417d-4b7d-9d2e-70884dde2daa.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a7a29e68ec614dc3a4d13cbe7cc940640000000/e2e/a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md", "", "", "a7a29e68-ec61-4dc3-a4d1-3cbe7cc94064.md") | Out-Null

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P5"))

Write-Host "Done applying handoff report update."
